# Apply the edits described by the diff against
# poc-professionnel/ig/StructureDefinition-CompetenceMetier.xlsx:
#   - Metadata!B5 (the "Title" row, previously blank) gets the value
#     "CompetenceMetier" (same text already used for the "Name" row, B4).
#   - Metadata!B8 (the "Date" row) is bumped from
#     2025-07-16T13:52:06+00:00 to 2025-07-17T14:35:50+00:00.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B5").Value = "CompetenceMetier"
$ws.Range("B8").Value = "2025-07-17T14:35:50+00:00"
